# Restore cell C10 on the "Rules" sheet from 18 to 1
# (per commit: Restored from revision #3f7ae98ccb1784fe27a553615da1d79a808b3cfb.TEST)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
